# Applies the "Updated cryptos list" data refresh to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Forces the cell to stay a text value (matches source t="inlineStr" cells)
    # even when $value looks like a number (e.g. "19.57"), without leaving any
    # residual cell style behind once we are done.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '25.787.11'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.633.54'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.43%  '
Set-TextValue 'D5' '215.61'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  +0.17%  '
Set-TextValue 'D9' '0.0638'
$ws.Range('E9').Value = '  -0.53%  '
Set-TextValue 'D10' '19.57'
$ws.Range('E10').Value = '  -3.01%  '
Set-TextValue 'D11' '0.0786'
$ws.Range('E11').Value = '  +0.80%  '
Set-TextValue 'D12' '4.25'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = '1.638.75'
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').Value = '1.862.43'
$ws.Range('E14').Value = '  +0.13%  '
Set-TextValue 'D15' '0.552'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').Value = '0.0₃0771'
$ws.Range('E16').Value = '  +1.16%  '
Set-TextValue 'D17' '63.34'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '25.813.09'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('E19').Value = '  +0.36%  '
Set-TextValue 'D20' '4.42'
$ws.Range('E20').Value = '  +1.76%  '
Set-TextValue 'D21' '193.38'
$ws.Range('E21').Value = '  -0.23%  '
Set-TextValue 'D22' '9.95'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  +1.42%  '
$ws.Range('E24').Value = '  +0.42%  '
Set-TextValue 'D25' '1.78'
$ws.Range('E25').Value = '  -0.33%  '
Set-TextValue 'D26' '139.64'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('E27').Value = '  -4.41%  '
Set-TextValue 'D28' '6.83'
$ws.Range('E28').Value = '  +0.26%  '
Set-TextValue 'D29' '15.58'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  -0.51%  '
Set-TextValue 'D32' '3.34'
$ws.Range('E32').Value = '  +0.93%  '
Set-TextValue 'D33' '3.26'
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('E35').Value = '  +0.77%  '
Set-TextValue 'D36' '0.895'
$ws.Range('E36').Value = '  -0.92%  '
Set-TextValue 'D37' '2.58'
$ws.Range('E37').Value = '  +0.40%  '
Set-TextValue 'D38' '0.549'
$ws.Range('E38').Value = '  -0.60%  '
$ws.Range('D39').Value = '1.106.93'
$ws.Range('E39').Value = '  -1.80%  '
Set-TextValue 'D40' '0.0156'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('E41').Value = '  +0.52%  '
Set-TextValue 'D42' '5.58'
$ws.Range('E42').Value = '  +1.43%  '
Set-TextValue 'D43' '99.69'
$ws.Range('E43').Value = '  +1.46%  '
Set-TextValue 'D44' '0.802'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '0.0₆0111'
$ws.Range('E45').Value = '  +0.52%  '
Set-TextValue 'D46' '55.11'
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D47' '0.419'
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '7.71'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.0504'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue 'D50' '2.34'
$ws.Range('E50').Value = '  +5.98%  '
$ws.Range('E51').Value = '  +0.56%  '
